# Update date strings in column A (rows 3-21) from "dd/mm/yyyy" to "dd-mm-yyyy"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Text
    if ($current -ne $null) {
        $updated = $current.Replace("/", "-")
        # Force text interpretation so Excel does not auto-convert the
        # dash-separated date string into a real date serial number.
        $cell.NumberFormat = "@"
        $cell.Value = $updated
        # Restore the cell's original (default) style now that the text
        # value is safely stored, so no visible formatting change remains.
        $cell.Style = "Normal"
    }
}
